$wb = $excel.ActiveWorkbook

# Capture the name of the existing (only) sheet before we add a new one.
$oldName = $wb.Worksheets.Item(1).Name

# Add a brand-new, unformatted worksheet (no custom column widths, no
# inherited cell styles) that will become the replacement "proveedores" sheet.
$ws = $wb.Worksheets.Add()
$ws.Name = "proveedores_new_tmp"

# Write the new (lowercase) header row into the fresh sheet.
$headers = @("codigo", "nombre", "nit", "telefono", "email", "direccion", "ciudad", "contacto", "comentario")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Remove the old "Datos" sheet (with its bold/blue header styling and the
# 20-wide custom columns) and rename the new sheet to take its place.
$old = $wb.Worksheets.Item($oldName)
$old.Delete()
$ws.Name = "proveedores"
